# Adds three new worksheets (CypherOutput_Message, StatOutput, StatOutput_Message)
# to the workbook, mirroring the "Message" sheet content and introducing a new
# aggregate "StatOutput" query/result pair, per commit "all canine test cases 72".

$wb = $excel.ActiveWorkbook

# ---- Source values, copied verbatim from the existing "Message" sheet ----
$neo4jUrlLabel   = "Neo4j_URL:"
$neo4jUrlValue   = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userNameLabel   = "User_name:"
$userNameValue   = "neo4j"
$pwdLabel        = "PWD:"
$pwdValue        = "icdcDBneo4j0"
$cypherLabel     = "Cypher:"
$cypherOutputQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_type IN ['Transcriptomics'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"
$outputLabel     = "Output:"
$outputPath      = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC02_Canine_Filter_StudyType-Transcriptomics_Neo4jData.xlsx"

# New aggregate ("stat") Cypher query introduced by this commit.
$statOutputQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_type IN ['Transcriptomics']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

function Add-MessageBlock($ws, $startRow, $cypherText) {
    $ws.Cells.Item($startRow,     1).Value = $neo4jUrlLabel
    $ws.Cells.Item($startRow + 1, 1).Value = $neo4jUrlValue
    $ws.Cells.Item($startRow + 2, 1).Value = $userNameLabel
    $ws.Cells.Item($startRow + 3, 1).Value = $userNameValue
    $ws.Cells.Item($startRow + 4, 1).Value = $pwdLabel
    $ws.Cells.Item($startRow + 5, 1).Value = $pwdValue
    $ws.Cells.Item($startRow + 6, 1).Value = $cypherLabel
    $ws.Cells.Item($startRow + 7, 1).Value = $cypherText
    $ws.Cells.Item($startRow + 8, 1).Value = $outputLabel
    $ws.Cells.Item($startRow + 9, 1).Value = $outputPath
}

# ---- Sheet 3: CypherOutput_Message (exact copy of Message) ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $afterSheet)
$ws3.Name = "CypherOutput_Message"
Add-MessageBlock $ws3 1 $cypherOutputQuery

# ---- Sheet 4: StatOutput (new aggregate counts table) ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $afterSheet)
$ws4.Name = "StatOutput"
$ws4.Cells.Item(1, 1).Value = "number_of_files"
$ws4.Cells.Item(1, 2).Value = "number_of_sample"
$ws4.Cells.Item(1, 3).Value = "number_of_cases"
$ws4.Cells.Item(1, 4).Value = "number_of_study"
$ws4.Cells.Item(2, 1).Value = "'331"
$ws4.Cells.Item(2, 2).Value = "'136"
$ws4.Cells.Item(2, 3).Value = "'60"
$ws4.Cells.Item(2, 4).Value = "'1"

# ---- Sheet 5: StatOutput_Message (Message block for CypherOutput query,
#      followed by a second Message block for the new StatOutput query) ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $afterSheet)
$ws5.Name = "StatOutput_Message"
Add-MessageBlock $ws5 1  $cypherOutputQuery
Add-MessageBlock $ws5 11 $statOutputQuery
